$p = $ppt.ActivePresentation

# --- Refresh the "datetimeFigureOut" date placeholder cached text
# (slide master + every slide layout) from 10/30/2024 to 11/5/2024,
# mirroring PowerPoint re-caching the auto date field on save.
$masterShapes = $p.SlideMaster.Shapes
for ($i = 1; $i -le $masterShapes.Count; $i++) {
  $shp = $masterShapes.Item($i)
  if ($shp.Name -like "Date Placeholder*") {
    $shp.TextFrame.TextRange.Text = "11/5/2024"
  }
}

$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
  $layoutShapes = $layouts.Item($i).Shapes
  for ($j = 1; $j -le $layoutShapes.Count; $j++) {
    $shp = $layoutShapes.Item($j)
    if ($shp.Name -like "Date Placeholder*") {
      $shp.TextFrame.TextRange.Text = "11/5/2024"
    }
  }
}

# --- Slide 9: "Work time: Lab 5" -> "Work time: Lab 4 or 5"
$slide = $p.Slides.Item(9)
$titleShape = $slide.Shapes.Item(1)
$tr = $titleShape.TextFrame.TextRange
$tr.Text = "Work time"
[void]$tr.InsertAfter(": Lab 4 or ")
[void]$tr.InsertAfter("5")
